$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Conclusion 1 paragraph: replace the "Roughly 50-60% ... successful.  "
# sentence with the expanded two-sentence version that also calls out the
# failure rate. Using Find/Replace on the exact old sentence keeps the
# preceding "Conclusion 1: " (bold) run and the single leading space run
# untouched, matching the rest of the paragraph's formatting (Times New
# Roman, no bold) for the newly written text.
# ---------------------------------------------------------------------------
$oldText = "Roughly 50-60% of crowdfunding projects in this dataset were considered successful.  "
$newText = "Of the campaigns that failed or succeeded, roughly 60% of crowdfunding projects in this dataset were considered successful and roughly 40% were considered to have failed."

$found = $d.Content.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)

if (-not $found) {
    throw "Could not find the Conclusion 1 sentence to replace."
}
